$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 20021
$ws.Range("I20").Value = 20021
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 20021
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -19791
$ws.Range("H35").Value = 20021
$ws.Range("I35").Value = 20021
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 20021
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -19642
$ws.Range("H69").Value = 1520.1515
$ws.Range("I69").Value = 2300
$ws.Range("J69").Value = 1495.7812
$ws.Range("K69").Value = 6900
$ws.Range("L69").Value = 4487.3436
$ws.Range("M69").Value = -6026
$ws.Range("N69").Value = -6235.3436
$ws.Range("H72").Value = 1520.1515
$ws.Range("I72").Value = 2300
$ws.Range("J72").Value = 1495.7812
$ws.Range("K72").Value = 20700
$ws.Range("L72").Value = 13462.0308
$ws.Range("M72").Value = -16332
$ws.Range("N72").Value = -22198.0308
$ws.Range("H121").Value = 1741.25
$ws.Range("I121").Value = 700
$ws.Range("J121").Value = 1949.5
$ws.Range("K121").Value = 2100
$ws.Range("L121").Value = 5848.5
$ws.Range("M121").Value = -353
$ws.Range("N121").Value = -9342.5
$ws.Range("H132").Value = 5473
$ws.Range("I132").Value = 6290.4546
$ws.Range("J132").Value = 3225
$ws.Range("K132").Value = 18871.3638
$ws.Range("L132").Value = 9675
$ws.Range("M132").Value = -16341.3638
$ws.Range("N132").Value = -14735
$ws.Range("H137").Value = 44841.87
$ws.Range("I137").Value = 1097.2307
$ws.Range("J137").Value = 101709.9
$ws.Range("K137").Value = 3291.6921
$ws.Range("L137").Value = 305129.7
$ws.Range("M137").Value = -741.6921000000002
$ws.Range("H138").Value = 1697.8793
$ws.Range("I138").Value = 655
$ws.Range("J138").Value = 2335.1943
$ws.Range("K138").Value = 1965
$ws.Range("L138").Value = 7005.5829
$ws.Range("M138").Value = 3175
$ws.Range("N138").Value = -17285.5829

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21113.445
$ws.Range("I32").Value = 21232.924
$ws.Range("J32").Value = 18007
$ws.Range("K32").Value = 21232.924
$ws.Range("L32").Value = 18007
$ws.Range("M32").Value = -20945.924
$ws.Range("N32").Value = -18581
$ws.Range("H61").Value = 2528.6191
$ws.Range("I61").Value = 1950.0555
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 1950.0555
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -1738.0555
$ws.Range("I63").Value = 900
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 900
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -214
$ws.Range("N63").Value = -31251372
$ws.Range("I66").Value = 900
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 4500
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -1068
$ws.Range("N66").Value = -156256864
$ws.Range("H74").Value = 45457492
$ws.Range("I74").Value = 58826496
$ws.Range("J74").Value = 2882.8
$ws.Range("K74").Value = 58826496
$ws.Range("L74").Value = 2882.8
$ws.Range("M74").Value = -58825622
$ws.Range("N74").Value = -4630.8
$ws.Range("H77").Value = 45457492
$ws.Range("I77").Value = 58826496
$ws.Range("J77").Value = 2882.8
$ws.Range("K77").Value = 294132480
$ws.Range("L77").Value = 14414
$ws.Range("M77").Value = -294128112
$ws.Range("N77").Value = -23150
$ws.Range("H101").Value = 40666.668
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 40666.668
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 40666.668
$ws.Range("N101").Value = -47156.668
$ws.Range("H136").Value = 2528.6191
$ws.Range("I136").Value = 1950.0555
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 5850.166499999999
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -3300.166499999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 5250
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 5250
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 5250
$ws.Range("N12").Value = -5586
$ws.Range("H82").Value = 12237.7
$ws.Range("I82").Value = 6315.778
$ws.Range("J82").Value = 65535
$ws.Range("K82").Value = 6315.778
$ws.Range("L82").Value = 65535
$ws.Range("M82").Value = -5932.778
$ws.Range("H85").Value = 12237.7
$ws.Range("I85").Value = 6315.778
$ws.Range("J85").Value = 65535
$ws.Range("K85").Value = 6315.778
$ws.Range("L85").Value = 65535
$ws.Range("M85").Value = -4989.778
$ws.Range("H88").Value = 14347.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 14347.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 14347.5
$ws.Range("N88").Value = -15159.5
$ws.Range("H91").Value = 14347.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 14347.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 14347.5
$ws.Range("N91").Value = -17155.5
$ws.Range("H94").Value = 834.96155
$ws.Range("I94").Value = 773.36365
$ws.Range("J94").Value = 1173.75
$ws.Range("K94").Value = 773.36365
$ws.Range("L94").Value = 1173.75
$ws.Range("M94").Value = -322.36365
$ws.Range("N94").Value = -2075.75
$ws.Range("H109").Value = 45000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 45000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 45000
$ws.Range("N109").Value = -47774
$ws.Range("H134").Value = 29049.64
$ws.Range("I134").Value = 38539.17
$ws.Range("J134").Value = 1530
$ws.Range("K134").Value = 115617.51
$ws.Range("L134").Value = 4590
$ws.Range("M134").Value = -113082.51
$ws.Range("N134").Value = -9660
$ws.Range("M88").ClearContents()
$ws.Range("M91").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13036.071
$ws.Range("I31").Value = 17285.55
$ws.Range("J31").Value = 2412.375
$ws.Range("K31").Value = 17285.55
$ws.Range("L31").Value = 2412.375
$ws.Range("M31").Value = -16990.55
$ws.Range("N31").Value = -3002.375
$ws.Range("H34").Value = 13036.071
$ws.Range("I34").Value = 17285.55
$ws.Range("J34").Value = 2412.375
$ws.Range("K34").Value = 17285.55
$ws.Range("L34").Value = 2412.375
$ws.Range("M34").Value = -17083.55
$ws.Range("N34").Value = -2816.375
$ws.Range("H94").Value = 4372.6665
$ws.Range("I94").Value = 2685.3333
$ws.Range("J94").Value = 6060
$ws.Range("K94").Value = 2685.3333
$ws.Range("L94").Value = 6060
$ws.Range("M94").Value = -2234.3333
$ws.Range("N94").Value = -6962
$ws.Range("H95").Value = 14000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 14000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 14000
$ws.Range("N95").Value = -19492
$ws.Range("H96").Value = 8307
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 8307
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 8307
$ws.Range("N96").Value = -13799
$ws.Range("H105").Value = 7353799
$ws.Range("I105").Value = 12500562
$ws.Range("J105").Value = 1280.1428
$ws.Range("K105").Value = 12500562
$ws.Range("L105").Value = 1280.1428
$ws.Range("M105").Value = -12498815
$ws.Range("N105").Value = -4774.1428
$ws.Range("H132").Value = 16377.833
$ws.Range("I132").Value = 17819.033
$ws.Range("J132").Value = 7442.4
$ws.Range("K132").Value = 53457.099
$ws.Range("L132").Value = 22327.2
$ws.Range("M132").Value = -50927.099
$ws.Range("N132").Value = -27387.2
$ws.Range("H134").Value = 981.6316
$ws.Range("I134").Value = 838.5714
$ws.Range("J134").Value = 1158.3529
$ws.Range("K134").Value = 2515.7142
$ws.Range("L134").Value = 3475.0587
$ws.Range("M134").Value = 19.28579999999965
$ws.Range("N134").Value = -8545.058700000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 296.18182
$ws.Range("I14").Value = 296.18182
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 888.54546
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -715.54546
$ws.Range("H131").Value = 702.54
$ws.Range("I131").Value = 590
$ws.Range("J131").Value = 706.0205999999999
$ws.Range("K131").Value = 1770
$ws.Range("L131").Value = 2118.0618
$ws.Range("M131").Value = 3270
$ws.Range("N131").Value = -12198.0618
$ws.Range("H132").Value = 1145.3636
$ws.Range("I132").Value = 828.4286
$ws.Range("J132").Value = 1700
$ws.Range("K132").Value = 7455.8574
$ws.Range("L132").Value = 15300
$ws.Range("M132").Value = -4925.8574
$ws.Range("N132").Value = -20360
$ws.Range("H139").Value = 1611.1578
$ws.Range("I139").Value = 1144.9333
$ws.Range("J139").Value = 3359.5
$ws.Range("K139").Value = 3434.7999
$ws.Range("L139").Value = 10078.5
$ws.Range("M139").Value = 1705.2001
$ws.Range("N139").Value = -20358.5
$ws.Range("H140").Value = 1315.8334
$ws.Range("I140").Value = 1071.875
$ws.Range("J140").Value = 3267.5
$ws.Range("K140").Value = 3215.625
$ws.Range("L140").Value = 9802.5
$ws.Range("M140").Value = 1964.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("H80").Value = 4200
$ws.Range("I80").Value = 3825
$ws.Range("J80").Value = 4450
$ws.Range("K80").Value = 3825
$ws.Range("L80").Value = 4450
$ws.Range("M80").Value = -2827
$ws.Range("N80").Value = -6446
$ws.Range("H83").Value = 4200
$ws.Range("I83").Value = 3825
$ws.Range("J83").Value = 4450
$ws.Range("K83").Value = 19125
$ws.Range("L83").Value = 22250
$ws.Range("M83").Value = -14133
$ws.Range("N83").Value = -32234
$ws.Range("H126").Value = 3668.725
$ws.Range("I126").Value = 2873.96
$ws.Range("J126").Value = 4993.3335
$ws.Range("K126").Value = 8621.880000000001
$ws.Range("L126").Value = 14980.0005
$ws.Range("M126").Value = -6151.880000000001
$ws.Range("N126").Value = -19920.0005
$ws.Range("M5").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 960.0714
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 890
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 890
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -1266

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1677.7778
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 1033.3334
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 2066.6668
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -4188.6668
$ws.Range("H84").Value = 1677.7778
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 1033.3334
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 10333.334
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -20941.334
$ws.Range("H103").Value = 15200.8
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 15200.8
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 15200.8
$ws.Range("N103").Value = -17544.8
